# Excel export stability fixes:
#  - rename/clarify date-of-birth header, add a new "WasDoo" (yes/no) column
#  - fix a typo in a sample birth-date value
#  - add a new sample data row
#  - normalize borders on the previously-blank trailing rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header row (row 1) ----
$ws.Range("E1").Value = "Дата рождения (ДД.ММ.ГГГГ)"
$ws.Range("F1").Value = "WasDoo (Да/Нет) или (1/0)"

# header formatting: gray fill, centered, new column F gets the same header look
$headerRange = $ws.Range("A1:F1")
$headerRange.Interior.Pattern = -4124      # xlSolid
$headerRange.Interior.PatternColorIndex = -4105
$headerRange.Interior.ThemeColor = 1
$headerRange.Interior.TintAndShade = -0.349986266670736
$headerRange.HorizontalAlignment = -4108   # xlCenter

$ws.Range("E1:F1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

# ---- fix existing typo in E3 ----
$ws.Range("E3").Value = "25.2.2000г"

# ---- new column F data ----
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = "нет"

# ---- new sample row 5 ----
$ws.Range("A5").Value = "test5"
$ws.Range("B5").Value = "test5"
$ws.Range("C5").Value = "test5"
$ws.Range("D5").Value = "1Г"
$ws.Range("E5").Value = "11.9.2001"
$ws.Range("F5").Value = "да"

# ---- borders ----
# full thin box border for all data cells A2:F12 (covers existing + new rows/col)
$dataRange = $ws.Range("A2:F12")
$dataRange.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$dataRange.Borders.Item(7).Weight = 2      # xlThin
$dataRange.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$dataRange.Borders.Item(8).Weight = 2
$dataRange.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$dataRange.Borders.Item(9).Weight = 2
$dataRange.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$dataRange.Borders.Item(10).Weight = 2
$dataRange.Borders.Item(11).LineStyle = 1  # xlInsideVertical
$dataRange.Borders.Item(11).Weight = 2
$dataRange.Borders.Item(12).LineStyle = 1  # xlInsideHorizontal
$dataRange.Borders.Item(12).Weight = 2

# header border: thin box around A1:F1 but no bottom edge (separation kept clean
# against the thin top border already drawn by the data range above)
$headerRange.Borders.Item(7).LineStyle = 1
$headerRange.Borders.Item(7).Weight = 2
$headerRange.Borders.Item(8).LineStyle = 1
$headerRange.Borders.Item(8).Weight = 2
$headerRange.Borders.Item(10).LineStyle = 1
$headerRange.Borders.Item(10).Weight = 2
$headerRange.Borders.Item(11).LineStyle = 1
$headerRange.Borders.Item(11).Weight = 2
$headerRange.Borders.Item(9).LineStyle = -4142 # xlLineStyleNone (bottom)

# ---- column width for the new column ----
$ws.Columns.Item(6).ColumnWidth = 17.5703125

# ---- selection ----
$ws.Range("A6").Select()
